$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "last updated" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 03:35"

# 2. Fix country order: swap Serbia (row 51) and Corea del Sur (row 52)
$ws.Range("A51").Value = "Corea del Sur"
$ws.Range("A52").Value = "Serbia"

# 3. Update statistics

# Estados Unidos (row 4)
$ws.Range("E4").Value = 1144734
$ws.Range("H4").Value = 100572

# Argentina (row 47)
$ws.Range("E47").Value = 8577
$ws.Range("H47").Value = 484

# Row 51 - now Corea del Sur (new/refreshed data)
$ws.Range("B51").Value = 11265
$ws.Range("C51").Value = 40
$ws.Range("D51").Value = 10295
$ws.Range("E51").Value = 701
$ws.Range("H51").Value = 269

# Row 52 - now Serbia (takes over previous Serbia numbers)
$ws.Range("B52").Value = 11227
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 6067
$ws.Range("E52").Value = 4921
$ws.Range("H52").Value = 239

# Australia (row 64)
$ws.Range("B64").Value = 7139
$ws.Range("C64").Value = 6
$ws.Range("D64").Value = 6560
$ws.Range("E64").Value = 476
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 103

# Jamaica (row 135)
$ws.Range("B135").Value = 564
$ws.Range("C135").Value = 8
$ws.Range("D135").Value = 267
$ws.Range("E135").Value = 288
